$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the value in A11 (was 4, now 5)
$ws.Range("A11").Value = 5

# Apply consistent 2-decimal number format to the data ranges
$ws.Range("A2:A13").NumberFormat = "0.00"
$ws.Range("B2:B13").NumberFormat = "0.00"

# Update the selected cell
$ws.Range("A12").Select() | Out-Null
